$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.258.52'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.35%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.647.82'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.76%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.38%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.03'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.46%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.08'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.65%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.35%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.586'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.14%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.647.41'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.75%  '

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.48%  '

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.48%  '

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.40%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.365'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.08%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.36'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.34%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.121.90'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.19%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.115.40'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.26%  '

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.76%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.662.26'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.63%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.41'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.77%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '342.92'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.27%  '

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.73%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.85'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.47%  '

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.23%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.10'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.59%  '

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.25%  '

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.90%  '

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.17%  '

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.24%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '546.53'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +15.40%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.08%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.90'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.71%  '

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.67%  '

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +6.82%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0810'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.28%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '172.36'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.19%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.10'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +11.85%  '

# Row 37
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.406'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.14%  '

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'FirstDigitalUSD'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.17%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.13'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.48%  '

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +6.51%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '172.03'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +8.70%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.16%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.76'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.30%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.38'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.18%  '

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +7.50%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.631'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.32%  '

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0962'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.08%  '

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0240'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.59%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.80'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.49%  '

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.12%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.23'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.29%  '
